$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 104
$ws.Range("F6").Value = 5138
$ws.Range("F7").Value = 420
$ws.Range("F9").Value = 894
$ws.Range("F17").Value = 1697
$ws.Range("F18").Value = 1436
$ws.Range("F19").Value = 785
$ws.Range("F22").Value = 288
$ws.Range("F27").Value = 523
$ws.Range("F28").Value = 2314
$ws.Range("F31").Value = 68
$ws.Range("F42").Value = 37
$ws.Range("F43").Value = 54

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 104
$ws.Range("F7").Value = 5138
$ws.Range("F8").Value = 420
$ws.Range("F12").Value = 894
$ws.Range("F23").Value = 1697
$ws.Range("F24").Value = 1436
$ws.Range("F25").Value = 785
$ws.Range("F28").Value = 288
$ws.Range("F33").Value = 523
$ws.Range("F34").Value = 2314
$ws.Range("F46").Value = 37
$ws.Range("F47").Value = 54
